$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- B7 (Experimental value): empty -> "false" ---
# NOTE: Assigning the literal text "false" directly to a Range.Value makes
# Excel auto-coerce it to the Boolean FALSE (same as typing it manually),
# which would change the cell's stored type/style. To keep it as a genuine
# text string (matching the target shared-string table), we build it as a
# text formula in a scratch cell, copy, and paste-special "Values" into the
# target cell - paste-special does not re-run literal-input type inference.
$scratch = $ws.Cells.Item(1, 5)
$scratch.Formula = "=""false"""
$scratch.Copy()
$ws.Cells.Item(7, 2).PasteSpecial(-4163)
$scratch.ClearContents()

# --- B8 (Date): old timestamp -> new timestamp (plain text) ---
$ws.Cells.Item(8, 2).Value = "2025-11-30T13:08:37+00:00"

# --- B17 (Description value): empty -> descriptive text ---
$ws.Cells.Item(17, 2).Value = "Units of measurement for VO2max values"
